$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 3).Value = 2.728779965065651
$ws.Cells.Item(2, 4).Value = 9.955051112981527
$ws.Cells.Item(2, 5).Value = 13.89494703328383
$ws.Cells.Item(2, 6).Value = 29.25823648110417
$ws.Cells.Item(2, 7).Value = 30.19113265778607
$ws.Cells.Item(2, 8).Value = 13.72557788070845
$ws.Cells.Item(2, 9).Value = 19.95349761985225
$ws.Cells.Item(2, 10).Value = 9.686377412883518
$ws.Cells.Item(2, 15).Value = 21.27759823122737

$ws.Cells.Item(3, 3).Value = 2.666799860447522
$ws.Cells.Item(3, 4).Value = 9.954949846502856
$ws.Cells.Item(3, 5).Value = 13.86112943064487
$ws.Cells.Item(3, 6).Value = 28.95767149847613
$ws.Cells.Item(3, 7).Value = 29.49286781126834
$ws.Cells.Item(3, 8).Value = 13.68904849605225
$ws.Cells.Item(3, 9).Value = 19.77431705042603
$ws.Cells.Item(3, 10).Value = 9.688278086992336
$ws.Cells.Item(3, 15).Value = 21.09485627290017

$ws.Cells.Item(4, 3).Value = 2.627431768047288
$ws.Cells.Item(4, 4).Value = 9.95657915190383
$ws.Cells.Item(4, 5).Value = 13.84336193971562
$ws.Cells.Item(4, 6).Value = 28.77966669272909
$ws.Cells.Item(4, 7).Value = 29.06527682824502
$ws.Cells.Item(4, 8).Value = 13.66951148769115
$ws.Cells.Item(4, 9).Value = 19.66886519832052
$ws.Cells.Item(4, 10).Value = 9.691131607468554
$ws.Cells.Item(4, 15).Value = 20.98749649895066

$ws.Cells.Item(5, 3).Value = 2.611069552182918
$ws.Cells.Item(5, 4).Value = 9.95766906031502
$ws.Cells.Item(5, 5).Value = 13.83687981959222
$ws.Cells.Item(5, 6).Value = 28.7088533455806
$ws.Cells.Item(5, 7).Value = 28.89160399499858
$ws.Cells.Item(5, 8).Value = 13.66228266299892
$ws.Cells.Item(5, 9).Value = 19.62708502643773
$ws.Cells.Item(5, 10).Value = 9.692718850239771
$ws.Cells.Item(5, 15).Value = 20.94500837609286

$ws.Cells.Item(6, 3).Value = 2.608333649158142
$ws.Cells.Item(6, 4).Value = 9.957875779935495
$ws.Cells.Item(6, 5).Value = 13.83584939071635
$ws.Cells.Item(6, 6).Value = 28.69720128737314
$ws.Cells.Item(6, 7).Value = 28.86280889360328
$ws.Cells.Item(6, 8).Value = 13.66112671683081
$ws.Cells.Item(6, 9).Value = 19.62022071728644
$ws.Cells.Item(6, 10).Value = 9.693008050396164
$ws.Cells.Item(6, 15).Value = 20.9380306942336

$ws.Cells.Item(7, 3).Value = 2.627212380296753
$ws.Cells.Item(7, 4).Value = 9.956592125435547
$ws.Cells.Item(7, 5).Value = 13.84327144362154
$ws.Cells.Item(7, 6).Value = 28.7787045920416
$ws.Cells.Item(7, 7).Value = 29.06293190951477
$ws.Cells.Item(7, 8).Value = 13.66941102414877
$ws.Cells.Item(7, 9).Value = 19.66829685213424
$ws.Cells.Item(7, 10).Value = 9.691151294850178
$ws.Cells.Item(7, 15).Value = 20.98691832433327

$ws.Cells.Item(8, 3).Value = 2.707687640182233
$ws.Cells.Item(8, 4).Value = 9.954665546593629
$ws.Cells.Item(8, 5).Value = 13.88266783122246
$ws.Cells.Item(8, 6).Value = 29.15329252787384
$ws.Cells.Item(8, 7).Value = 29.95030216823793
$ws.Cells.Item(8, 8).Value = 13.71238573299202
$ws.Cells.Item(8, 9).Value = 19.890795804572
$ws.Cells.Item(8, 10).Value = 9.686682972881327
$ws.Cells.Item(8, 15).Value = 21.21361027950487

$ws.Cells.Item(9, 3).Value = 2.854662956744724
$ws.Cells.Item(9, 4).Value = 9.96427567094139
$ws.Cells.Item(9, 5).Value = 13.98346104558504
$ws.Cells.Item(9, 6).Value = 29.93610847079733
$ws.Cells.Item(9, 7).Value = 31.68738879936633
$ws.Cells.Item(9, 8).Value = 13.81933832538164
$ws.Cells.Item(9, 9).Value = 20.36128806603438
$ws.Cells.Item(9, 10).Value = 9.691282453223312
$ws.Cells.Item(9, 15).Value = 21.69455720922446

$ws.Cells.Item(10, 3).Value = 2.955520451787618
$ws.Cells.Item(10, 4).Value = 9.979444865072413
$ws.Cells.Item(10, 5).Value = 14.07149579450333
$ws.Cells.Item(10, 6).Value = 30.53555718835969
$ws.Cells.Item(10, 7).Value = 32.94627397956596
$ws.Cells.Item(10, 8).Value = 13.91133742855725
$ws.Cells.Item(10, 9).Value = 20.72490527758731
$ws.Cells.Item(10, 10).Value = 9.702772822597673
$ws.Cells.Item(10, 15).Value = 22.06723566228228

$ws.Cells.Item(11, 3).Value = 2.999759994391921
$ws.Cells.Item(11, 4).Value = 9.988091451713009
$ws.Cells.Item(11, 5).Value = 14.11448837401452
$ws.Cells.Item(11, 6).Value = 30.8124209800499
$ws.Cells.Item(11, 7).Value = 33.5121897987804
$ws.Cells.Item(11, 8).Value = 13.95600390532065
$ws.Cells.Item(11, 9).Value = 20.89357776281338
$ws.Cells.Item(11, 10).Value = 9.709750526752767
$ws.Cells.Item(11, 15).Value = 22.24033131613716

$ws.Cells.Item(12, 3).Value = 3.016269201105552
$ws.Cells.Item(12, 4).Value = 9.991615157128905
$ws.Cells.Item(12, 5).Value = 14.1311831220988
$ws.Cells.Item(12, 6).Value = 30.91776706014752
$ws.Cells.Item(12, 7).Value = 33.72528506671544
$ws.Cells.Item(12, 8).Value = 13.97331327431777
$ws.Cells.Item(12, 9).Value = 20.95786339543751
$ws.Cells.Item(12, 10).Value = 9.712643302057845
$ws.Cells.Item(12, 15).Value = 22.30633505977239

$ws.Cells.Item(13, 3).Value = 3.012724581529501
$ws.Cells.Item(13, 4).Value = 9.990845198871439
$ws.Cells.Item(13, 5).Value = 14.12756933111613
$ws.Cells.Item(13, 6).Value = 30.89505801903287
$ws.Cells.Item(13, 7).Value = 33.67944830912072
$ws.Cells.Item(13, 8).Value = 13.9695679794355
$ws.Cells.Item(13, 9).Value = 20.94400085894352
$ws.Cells.Item(13, 10).Value = 9.712009173623047
$ws.Cells.Item(13, 15).Value = 22.29210058217862

$ws.Cells.Item(14, 3).Value = 3.001123133745226
$ws.Cells.Item(14, 4).Value = 9.988376360225647
$ws.Cells.Item(14, 5).Value = 14.11585361043715
$ws.Cells.Item(14, 6).Value = 30.82107831877408
$ws.Cells.Item(14, 7).Value = 33.52974669752562
$ws.Cells.Item(14, 8).Value = 13.95742009621922
$ws.Cells.Item(14, 9).Value = 20.89885862447361
$ws.Cells.Item(14, 10).Value = 9.709983505264406
$ws.Cells.Item(14, 15).Value = 22.24575266440347

$ws.Cells.Item(15, 3).Value = 2.99398501318501
$ws.Cells.Item(15, 4).Value = 9.98689655637704
$ws.Cells.Item(15, 5).Value = 14.10873107781515
$ws.Cells.Item(15, 6).Value = 30.77582637344957
$ws.Cells.Item(15, 7).Value = 33.43788660025978
$ws.Cells.Item(15, 8).Value = 13.9500303268308
$ws.Cells.Item(15, 9).Value = 20.87125982419021
$ws.Cells.Item(15, 10).Value = 9.708775302607123
$ws.Cells.Item(15, 15).Value = 22.21742096088007

$ws.Cells.Item(16, 3).Value = 2.952595627115497
$ws.Cells.Item(16, 4).Value = 9.978914788576947
$ws.Cells.Item(16, 5).Value = 14.0687446273984
$ws.Cells.Item(16, 6).Value = 30.51753921407939
$ws.Cells.Item(16, 7).Value = 32.90913378270869
$ws.Cells.Item(16, 8).Value = 13.90847420835361
$ws.Cells.Item(16, 9).Value = 20.71394298549209
$ws.Cells.Item(16, 10).Value = 9.702351936432365
$ws.Cells.Item(16, 15).Value = 22.0559903562594

$ws.Cells.Item(17, 3).Value = 2.926778821510534
$ws.Cells.Item(17, 4).Value = 9.974464323503957
$ws.Cells.Item(17, 5).Value = 14.04496227410565
$ws.Cells.Item(17, 6).Value = 30.3600884634232
$ws.Cells.Item(17, 7).Value = 32.58285940051093
$ws.Cells.Item(17, 8).Value = 13.8836952925915
$ws.Cells.Item(17, 9).Value = 20.61822981597839
$ws.Cells.Item(17, 10).Value = 9.698858949691044
$ws.Cells.Item(17, 15).Value = 21.95783059923031

$ws.Cells.Item(18, 3).Value = 2.911775614143854
$ws.Cells.Item(18, 4).Value = 9.972068990275975
$ws.Cells.Item(18, 5).Value = 14.03156109064438
$ws.Cells.Item(18, 6).Value = 30.26992536866278
$ws.Cells.Item(18, 7).Value = 32.39457190214151
$ws.Cells.Item(18, 8).Value = 13.86970875867027
$ws.Cells.Item(18, 9).Value = 20.56348847637523
$ws.Cells.Item(18, 10).Value = 9.697014745871311
$ws.Cells.Item(18, 15).Value = 21.90171041309423

$ws.Cells.Item(19, 3).Value = 2.906669560474804
$ws.Cells.Item(19, 4).Value = 9.971286264311502
$ws.Cells.Item(19, 5).Value = 14.02707165926478
$ws.Cells.Item(19, 6).Value = 30.23946905090553
$ws.Cells.Item(19, 7).Value = 32.33072133737632
$ws.Cells.Item(19, 8).Value = 13.86501906968423
$ws.Cells.Item(19, 9).Value = 20.54500899855126
$ws.Cells.Item(19, 10).Value = 9.696418682596818
$ws.Cells.Item(19, 15).Value = 21.88276895557813

$ws.Cells.Item(20, 3).Value = 2.929543064964538
$ws.Cells.Item(20, 4).Value = 9.974921073813178
$ws.Cells.Item(20, 5).Value = 14.04746525924699
$ws.Cells.Item(20, 6).Value = 30.37680882686303
$ws.Cells.Item(20, 7).Value = 32.61765809388438
$ws.Cells.Item(20, 8).Value = 13.88630562597058
$ws.Cells.Item(20, 9).Value = 20.62838693071233
$ws.Cells.Item(20, 10).Value = 9.69921372935608
$ws.Cells.Item(20, 15).Value = 21.96824522223619

$ws.Cells.Item(21, 3).Value = 3.004537420118909
$ws.Cells.Item(21, 4).Value = 9.989094762572616
$ws.Cells.Item(21, 5).Value = 14.11928363092765
$ws.Cells.Item(21, 6).Value = 30.84279504255866
$ws.Cells.Item(21, 7).Value = 33.57375216391741
$ws.Cells.Item(21, 8).Value = 13.96097758055037
$ws.Cells.Item(21, 9).Value = 20.91210723109233
$ws.Cells.Item(21, 10).Value = 9.710571706527887
$ws.Cells.Item(21, 15).Value = 22.25935424026518

$ws.Cells.Item(22, 3).Value = 3.052129177597081
$ws.Cells.Item(22, 4).Value = 9.999811063477841
$ws.Cells.Item(22, 5).Value = 14.16863180903648
$ws.Cells.Item(22, 6).Value = 31.15023604654399
$ws.Cells.Item(22, 7).Value = 34.19149596608563
$ws.Cells.Item(22, 8).Value = 14.01207827847209
$ws.Cells.Item(22, 9).Value = 21.09991615645757
$ws.Cells.Item(22, 10).Value = 9.719454015011983
$ws.Cells.Item(22, 15).Value = 22.45224346823443

$ws.Cells.Item(23, 3).Value = 3.026860887910352
$ws.Cells.Item(23, 4).Value = 9.993959217700572
$ws.Cells.Item(23, 5).Value = 14.14207635031708
$ws.Cells.Item(23, 6).Value = 30.98591603219004
$ws.Cells.Item(23, 7).Value = 33.86251865599009
$ws.Cells.Item(23, 8).Value = 13.98459795704715
$ws.Cells.Item(23, 9).Value = 20.99947965715831
$ws.Cells.Item(23, 10).Value = 9.714580309249527
$ws.Cells.Item(23, 15).Value = 22.34907260818603

$ws.Cells.Item(24, 3).Value = 2.928293850742838
$ws.Cells.Item(24, 4).Value = 9.974714068218672
$ws.Cells.Item(24, 5).Value = 14.04633281295985
$ws.Cells.Item(24, 6).Value = 30.36924842973588
$ws.Cells.Item(24, 7).Value = 32.60192779189889
$ws.Cells.Item(24, 8).Value = 13.8851246860798
$ws.Cells.Item(24, 9).Value = 20.62379400570167
$ws.Cells.Item(24, 10).Value = 9.699052822652378
$ws.Cells.Item(24, 15).Value = 21.96353579133698

$ws.Cells.Item(25, 3).Value = 2.816118648599054
$ws.Cells.Item(25, 4).Value = 9.960247339920496
$ws.Cells.Item(25, 5).Value = 13.9537090300142
$ws.Cells.Item(25, 6).Value = 29.71968912420404
$ws.Cells.Item(25, 7).Value = 31.21943570676752
$ws.Cells.Item(25, 8).Value = 13.78801584276745
$ws.Cells.Item(25, 9).Value = 20.23064075146674
$ws.Cells.Item(25, 10).Value = 9.688610579383473
$ws.Cells.Item(25, 15).Value = 21.56084013990481
